$wb = $excel.ActiveWorkbook

# Work on the "week2" sheet (the active tab), where this week's DQ1 time
# tracking lives.
$ws = $wb.Worksheets.Item("week2")

# C2 ("Discussion question 1" row): actual time to complete updated from
# 1:30 (0.0625) to 2:50 (0.11805555555555557).
$ws.Range("C2").Value = 0.11805555555555557

# C3 ("Discussion question 2" row): actual time now logged as 1:15
# (0.052083333333333336); previously blank.
$ws.Range("C3").Value = 0.052083333333333336

# C19 total: extend the SUM range to include row 2 (it previously started
# at row 4, skipping the newly-filled C2/C3 values).
$ws.Range("C19").Formula = "=SUM(C2:C18)"

# Reposition the saved window (cosmetic, matches the author's last on-screen
# window placement when the file was saved).
$excel.ActiveWindow.Left = 9020
$excel.ActiveWindow.Top = 2560
